# =====================================================================
# Plan ECO-FINAN-PUM.xlsx -- "calculados indicadores financieros para RE"
#
# Adds three new worksheets (Ganancias, Costos, Indicadores Financieros)
# between "Ingresos" and "Precios ", wires them up with cross-sheet
# formulas, tweaks a handful of input cells on "Precios " / "Volumenes",
# and adds a small yearly-total table to "Ingresos".
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the three new sheets, in order, right after "Ingresos".
# ---------------------------------------------------------------------
$wsIngresos = $wb.Worksheets.Item("Ingresos")

$wsGanancias = $wb.Worksheets.Add($null, $wsIngresos)
$wsGanancias.Name = "Ganancias"

$wsCostos = $wb.Worksheets.Add($null, $wsGanancias)
$wsCostos.Name = "Costos"

$wsIndicadores = $wb.Worksheets.Add($null, $wsCostos)
$wsIndicadores.Name = "Indicadores Financieros"

$wsPrecios = $wb.Worksheets.Item("Precios ")
$wsVolumenes = $wb.Worksheets.Item("Volumenes")

# ---------------------------------------------------------------------
# 2) "Ingresos" -- add the yearly totals table (rows 10-11).
# ---------------------------------------------------------------------
$wsIngresos.Range("B10").Value2 = "Año"
$wsIngresos.Range("C10").Value2 = 2016
$wsIngresos.Range("D10").Value2 = 2017

$wsIngresos.Range("B11").Value2 = "Ingreso"
$wsIngresos.Range("C11").Formula = "=SUM(F3:F4)"
$wsIngresos.Range("D11").Formula = "=SUM(F5:F6)"

# ---------------------------------------------------------------------
# 3) "Precios " -- update base prices (recalculates B14/B15/B17 too).
# ---------------------------------------------------------------------
$wsPrecios.Range("B5").Value2 = 1700
$wsPrecios.Range("B6").Value2 = 6500
$wsPrecios.Range("B8").Value2 = 1400

# ---------------------------------------------------------------------
# 4) "Volumenes" -- update volumes for semesters 2016-2 and the whole
#    "Subscripciones estándar/premium Negocios" rows.
# ---------------------------------------------------------------------
$wsVolumenes.Range("F4").Value2 = 4
$wsVolumenes.Range("G4").Value2 = 5
$wsVolumenes.Range("H4").Value2 = 5

$wsVolumenes.Range("E5").Value2 = 5
$wsVolumenes.Range("H5").Value2 = 7

$wsVolumenes.Range("B6").Value2 = 13
$wsVolumenes.Range("C6").Value2 = 7
$wsVolumenes.Range("D6").Value2 = 4000
$wsVolumenes.Range("E6").Value2 = 6
$wsVolumenes.Range("F6").Value2 = 5
$wsVolumenes.Range("G6").Value2 = 5
$wsVolumenes.Range("H6").Value2 = 8

$wsVolumenes.Range("B7").Value2 = 21
$wsVolumenes.Range("C7").Value2 = 9
$wsVolumenes.Range("D7").Value2 = 2000
$wsVolumenes.Range("E7").Value2 = 5
$wsVolumenes.Range("F7").Value2 = 3
$wsVolumenes.Range("G7").Value2 = 6
$wsVolumenes.Range("H7").Value2 = 9

# ---------------------------------------------------------------------
# 5) "Ganancias" -- Ingresos - Egresos - impuesto a las ganancias.
# ---------------------------------------------------------------------
$wsGanancias.Range("B1").Value2 = "Año 1"
$wsGanancias.Range("C1").Value2 = "Año 2"

$wsGanancias.Range("A2").Value2 = "Ingresos"
$wsGanancias.Range("B2").Formula = "=Ingresos!C11"
$wsGanancias.Range("C2").Formula = "=Ingresos!D11"

$wsGanancias.Range("A3").Value2 = "Egresos"
$wsGanancias.Range("B3").Formula = "=Costos!B5"
$wsGanancias.Range("C3").Formula = "=Costos!C5"

$wsGanancias.Range("G3").Value2 = "GANANCIAS"
$wsGanancias.Range("H3").Value2 = 0.35
$wsGanancias.Range("H3").NumberFormat = "0%"

$wsGanancias.Range("A4").Value2 = "Ganancia Neta"
$wsGanancias.Range("B4").Formula = "=B2-B3"
$wsGanancias.Range("C4").Formula = "=C2-C3"

$wsGanancias.Range("A5").Value2 = "imp ganancia"
$wsGanancias.Range("B5").Formula = "=-(B4*H3)"
$wsGanancias.Range("C5").Formula = "=-(C4*H3)"

$wsGanancias.Range("A6").Value2 = "total"
$wsGanancias.Range("B6").Formula = "=SUM(B4:B5)"
$wsGanancias.Range("C6").Formula = "=SUM(C4:C5)"

$wsGanancias.Columns.Item(1).ColumnWidth = 13.42578125
$wsGanancias.Columns.Item(2).ColumnWidth = 11.85546875

# ---------------------------------------------------------------------
# 6) "Costos" -- gastos fijos / variables / aguinaldos.
# ---------------------------------------------------------------------
$wsCostos.Range("B1").Value2 = "Año 1"
$wsCostos.Range("C1").Value2 = "Año 2"

$wsCostos.Range("A2").Value2 = "Gastos Fijos"
$wsCostos.Range("B2").Value2 = 2086920
$wsCostos.Range("C2").Formula = "=B2+(B2*'Precios '!G3)+503577"

$wsCostos.Range("A3").Value2 = "Gastos Variables"
$wsCostos.Range("B3").Value2 = 420000
$wsCostos.Range("C3").Formula = "=B3+(B3*'Precios '!G4)+128250"

$wsCostos.Range("A4").Value2 = "Aguinaldos"
$wsCostos.Range("B4").Value2 = 122168
$wsCostos.Range("C4").Value2 = 171035

$wsCostos.Range("A5").Value2 = "total"
$wsCostos.Range("B5").Formula = "=SUM(B2:B4)"
$wsCostos.Range("C5").Formula = "=SUM(C2:C4)"
$wsCostos.Range("B5:C5").Interior.Color = 65535

$wsCostos.Columns.Item(1).ColumnWidth = 15.5703125
$wsCostos.Columns.Item(2).ColumnWidth = 15.5703125

# ---------------------------------------------------------------------
# 7) "Indicadores Financieros" -- VAN / TIR / payback.
# ---------------------------------------------------------------------
$wsIndicadores.Range("A2").Value2 = "año"
$wsIndicadores.Range("B2").Value2 = 0
$wsIndicadores.Range("C2").Value2 = 1
$wsIndicadores.Range("D2").Value2 = 2

$wsIndicadores.Range("B3").Value2 = -342853
$wsIndicadores.Range("C3").Formula = "=Ganancias!B6"
$wsIndicadores.Range("D3").Formula = "=Ganancias!C6"

$wsIndicadores.Range("H4").Value2 = "TASA"
$wsIndicadores.Range("I4").Value2 = 0.234
$wsIndicadores.Range("I4").NumberFormat = "0.00%"

$wsIndicadores.Range("A9").Value2 = "VAN"
$wsIndicadores.Range("B9").Formula = "=NPV(I4,C3:D3)+B3"
$wsIndicadores.Range("B9").NumberFormat = """$""\ #,##0.00;[Red]""$""\ \-#,##0.00"

$wsIndicadores.Range("A10").Value2 = "TIR"
$wsIndicadores.Range("B10").Formula = "=IRR(B3:D3)"
$wsIndicadores.Range("B10").NumberFormat = "0.00%"

$wsIndicadores.Range("A13").Value2 = "año"
$wsIndicadores.Range("B13").Value2 = 0
$wsIndicadores.Range("C13").Value2 = 1
$wsIndicadores.Range("D13").Value2 = 2

$wsIndicadores.Range("A14").Value2 = "cashflow"
$wsIndicadores.Range("B14").Formula = "=B3"
$wsIndicadores.Range("C14").Formula = "=B14+C3"
$wsIndicadores.Range("D14").Formula = "=C14+D3"

$wsIndicadores.Range("C16").Value2 = "meses"

$wsIndicadores.Range("B17").Formula = "=D3"
$wsIndicadores.Range("C17").Value2 = 12

$wsIndicadores.Range("B18").Formula = "=-(C14)"
$wsIndicadores.Range("C18").Formula = "=(B18*C17)/B17"

$wsIndicadores.Range("B20").Value2 = "1 año"
$wsIndicadores.Range("C20").Value2 = "3 meses"
$wsIndicadores.Range("D20").Value2 = "6 dias"

$wsIndicadores.Columns.Item(2).ColumnWidth = 13.140625
$wsIndicadores.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 8) Selections -- match each sheet's final cursor position, then make
#    "Indicadores Financieros" the active tab (author's last stop).
# ---------------------------------------------------------------------
$wsIngresos.Range("C11").Select()
$wsGanancias.Range("B6").Select()
$wsCostos.Range("D11").Select()
$wsPrecios.Range("B4").Select()
$wsVolumenes.Range("B7").Select()

$wsIndicadores.Activate()
$wsIndicadores.Range("C21").Select()

Write-Output "done"
